# ROADMAP.docx edit script
# Applies the Lithuanian-translations commit: swaps several roadmap item
# headings, drops five stale "Completed Features"-adjacent headings,
# restyles/empties a couple of paragraphs and rewrites + splits the long
# "Completed Features" summary paragraph, finishing with a new blank
# styled paragraph before "Community Projects".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Six simple heading text swaps (Find/Replace scoped to the exact
#    paragraph so we never touch a look-alike elsewhere in the doc).
# ---------------------------------------------------------------------
function Replace-ParagraphText($paraIndex, $oldText, $newText) {
    $rng = $d.Paragraphs.Item($paraIndex).Range
    $null = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

Replace-ParagraphText 4 "SMARTCARD AND MERCHANT READER" "SMARTNODE STARTING FROM WEBWALLET"
Replace-ParagraphText 5 "HARDWARE WALLETS" "NATIVE MOBILE WALLET WITH MULTICURRENCY TRADING ABILITY"
Replace-ParagraphText 6 "UPGRADED EXPLORER" "FUNCTION FOR TRANSACTION LIMITED MINING"
Replace-ParagraphText 7 "SMARTNODE STARTING FROM WEBWALLET" "SUPERNODES"
Replace-ParagraphText 8 "NATIVE MOBILE WALLET WITH MULTICURRENCY TRADING ABILITY" "IMPROVED SMARTREWARDS"
Replace-ParagraphText 9 "ADAPTIVE BLOCKS" "ON BLOCKCHAIN PROPOSAL VOTING"

# ---------------------------------------------------------------------
# 2) Delete the five now-duplicate/stale heading paragraphs:
#    24HR SUPPORT CENTER, SUPERNODES, IMPROVED SMARTREWARDS,
#    SMARTNODE STARTING FROM ELECTRUM, ON BLOCKCHAIN PROPOSAL VOTING
#    (paragraphs 10..14, right before "COLD VOTING KEYS").
# ---------------------------------------------------------------------
$startDel = $d.Paragraphs.Item(10).Range.Start
$endDel = $d.Paragraphs.Item(14).Range.End
$d.Range($startDel, $endDel).Delete()

# ---------------------------------------------------------------------
# 3) "COLD VOTING KEYS" paragraph: add shading, drop the after-spacing to
#    0 and switch the (inherit/Times New Roman) fonts to "open sans".
# ---------------------------------------------------------------------
$coldVotingXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="396" w:lineRule="atLeast"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="open sans" w:eastAsia="Times New Roman" w:hAnsi="open sans" w:cs="open sans"/><w:b/><w:bCs/><w:caps/><w:color w:val="343434"/><w:spacing w:val="8"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="open sans" w:hAnsi="open sans" w:cs="open sans" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:caps/><w:color w:val="343434"/><w:spacing w:val="8"/></w:rPr><w:t xml:space="preserve">COLD VOTING KEYS</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$coldVotingPara = $d.Paragraphs.Item(10)
$coldVotingPara.Range.InsertXML($coldVotingXml)

# ---------------------------------------------------------------------
# 4) "ELECTRUM SMARTVOTING" paragraph loses its run entirely, leaving an
#    empty paragraph with the same paragraph properties.
# ---------------------------------------------------------------------
$electrumXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="396" w:lineRule="atLeast"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="open sans" w:eastAsia="Times New Roman" w:hAnsi="open sans" w:cs="open sans"/><w:b/><w:bCs/><w:caps/><w:color w:val="343434"/><w:spacing w:val="8"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$electrumPara = $d.Paragraphs.Item(12)
$electrumPara.Range.InsertXML($electrumXml)

# ---------------------------------------------------------------------
# 5) The long "Completed Features" summary paragraph is rewritten and
#    split into three paragraphs:
#      a) empty paragraph carrying the old run's rPr as its pPr/rPr
#      b) the summary paragraph itself (old pPr), with new text
#      c) a new empty "section-header style" paragraph (Tahoma/F4B517)
# ---------------------------------------------------------------------
$summaryXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="open sans" w:eastAsia="Times New Roman" w:hAnsi="open sans" w:cs="open sans"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="open sans" w:hAnsi="open sans" w:cs="open sans" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">SmartRewards, SmartHive voting, SmartHive, Smart Web wallet , Mobile wallets, SmartNodes, Electrum, Pay to email, InstantPay, Core upgrades for faster syncing, vault, multisig (complete), trezor (hardware wallet), upgraded explorer, 24hr support center, smartnode starting from electrum, electrum smartvoting, SmartCard and Merchant Reader, Collateral change to 100k to enable better quality SmartNodes, SmartShift, SmartRewards tab in Node and Electrum wallets, SmartNode starting with Trezor, SmartCard and Merchant Reader</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="264" w:lineRule="atLeast"/><w:textAlignment w:val="baseline"/><w:outlineLvl w:val="2"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="Times New Roman" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:color w:val="F4B517"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$summaryPara = $d.Paragraphs.Item(14)
$summaryPara.Range.InsertXML($summaryXml)
